$p = $ppt.ActivePresentation

# 1. Update the cached "datetimeFigureOut" date text (1/1/2016 -> 1/2/2016)
#    on the slide master and on every slide layout's Date placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "1/1/2016") {
        $sh.TextFrame.TextRange.Text = "1/2/2016"
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "1/1/2016") {
            $sh.TextFrame.TextRange.Text = "1/2/2016"
        }
    }
}

# 2. Change the "temp <= T-LO" label on slide 1 to "temp < T-LO", splitting
#    the "<=" glyph + trailing space into its own run (< ) while keeping the
#    rest of the text in separate runs, matching the authored edit.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "temp $([char]0x2264) T-LO") {
        $tr = $sh.TextFrame.TextRange
        $chars = $tr.Characters(6, 2)
        $chars.Text = "< "
    }
}
